$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the Price (D) and Volume(1h) (E) columns for each coin row.
# Values in column D that look numeric must be written with a leading
# apostrophe so Excel stores them as literal text (preserving formats
# such as "1.00" or "47.236.20") instead of converting them to numbers.

$ws.Range("D2").Value = '47.236.20'
$ws.Range("E2").Value = '  +4.34%  '
$ws.Range("D3").Value = '2.486.19'
$ws.Range("E3").Value = '  +2.04%  '
$ws.Range("D4").Value = "'" + '0.999'
$ws.Range("E4").Value = '  -0.08%  '
$ws.Range("D5").Value = "'" + '322.35'
$ws.Range("E5").Value = '  +1.20%  '
$ws.Range("D6").Value = "'" + '104.99'
$ws.Range("E6").Value = '  +1.19%  '
$ws.Range("D7").Value = "'" + '0.521'
$ws.Range("E7").Value = '  +1.08%  '
$ws.Range("D8").Value = "'" + '1.00'
$ws.Range("E8").Value = '  +0.05%  '
$ws.Range("D9").Value = "'" + '0.538'
$ws.Range("E9").Value = '  +1.82%  '
$ws.Range("D10").Value = "'" + '37.28'
$ws.Range("E10").Value = '  +4.03%  '
$ws.Range("D11").Value = "'" + '0.0811'
$ws.Range("E11").Value = '  +0.97%  '
$ws.Range("D12").Value = "'" + '0.123'
$ws.Range("E12").Value = '  +0.07%  '
$ws.Range("D13").Value = "'" + '18.25'
$ws.Range("E13").Value = '  -0.83%  '
$ws.Range("D14").Value = "'" + '7.16'
$ws.Range("E14").Value = '  +2.42%  '
$ws.Range("D15").Value = '2.874.10'
$ws.Range("E15").Value = '  +1.94%  '
$ws.Range("D16").Value = '2.498.50'
$ws.Range("E16").Value = '  +2.60%  '
$ws.Range("D17").Value = "'" + '0.841'
$ws.Range("E17").Value = '  +1.14%  '
$ws.Range("D18").Value = '47.148.51'
$ws.Range("E18").Value = '  +4.41%  '
$ws.Range("D19").Value = "'" + '12.70'
$ws.Range("E19").Value = '  +3.31%  '
$ws.Range("D20").Value = "'" + '6.53'
$ws.Range("E20").Value = '  +2.31%  '
$ws.Range("D21").Value = '0.0₃0932'
$ws.Range("E21").Value = '  +0.63%  '
$ws.Range("D22").Value = "'" + '70.55'
$ws.Range("E22").Value = '  +2.15%  '
$ws.Range("D23").Value = "'" + '250.01'
$ws.Range("E23").Value = '  +2.69%  '
$ws.Range("D24").Value = "'" + '2.36'
$ws.Range("E24").Value = '  +3.43%  '
$ws.Range("D25").Value = "'" + '2.55'
$ws.Range("E25").Value = '  +1.47%  '
$ws.Range("D26").Value = "'" + '26.10'
$ws.Range("E26").Value = '  +2.70%  '
$ws.Range("D27").Value = "'" + '1.00'
$ws.Range("E27").Value = '  -0.05%  '
$ws.Range("D28").Value = "'" + '10.09'
$ws.Range("E28").Value = '  +5.64%  '
$ws.Range("D29").Value = "'" + '2.20'
$ws.Range("E29").Value = '  +0.77%  '
$ws.Range("D30").Value = "'" + '35.32'
$ws.Range("E30").Value = '  +5.66%  '
$ws.Range("D31").Value = "'" + '0.134'
$ws.Range("E31").Value = '  +4.82%  '
$ws.Range("D32").Value = "'" + '49.49'
$ws.Range("E32").Value = '  +0.24%  '
$ws.Range("D33").Value = "'" + '19.79'
$ws.Range("E33").Value = '  -3.22%  '
$ws.Range("D34").Value = "'" + '5.37'
$ws.Range("E34").Value = '  +2.78%  '
$ws.Range("D35").Value = "'" + '0.0780'
$ws.Range("E35").Value = '  +1.66%  '
$ws.Range("D36").Value = "'" + '1.00'
$ws.Range("E36").Value = '  +0.12%  '
$ws.Range("D37").Value = "'" + '4.62'
$ws.Range("E37").Value = '  +2.73%  '
$ws.Range("D38").Value = "'" + '1.93'
$ws.Range("E38").Value = '  +1.70%  '
$ws.Range("D39").Value = "'" + '2.96'
$ws.Range("E39").Value = '  +3.71%  '
$ws.Range("D40").Value = "'" + '0.111'
$ws.Range("E40").Value = '  +1.25%  '
$ws.Range("D41").Value = "'" + '121.51'
$ws.Range("E41").Value = '  -2.32%  '
$ws.Range("D42").Value = "'" + '2.22'
$ws.Range("E42").Value = '  +0.90%  '
$ws.Range("D43").Value = "'" + '21.36'
$ws.Range("E43").Value = '  -0.39%  '
$ws.Range("D44").Value = "'" + '0.0294'
$ws.Range("E44").Value = '  +1.23%  '
$ws.Range("D45").Value = '1.950.77'
$ws.Range("E45").Value = '  +0.48%  '
$ws.Range("D46").Value = "'" + '2.97'
$ws.Range("E46").Value = '  +1.14%  '
$ws.Range("D47").Value = "'" + '2.10'
$ws.Range("E47").Value = '  +0.19%  '
$ws.Range("D48").Value = "'" + '9.19'
$ws.Range("E48").Value = '  -0.69%  '
$ws.Range("D49").Value = "'" + '1.79'
$ws.Range("E49").Value = '  +0.89%  '
$ws.Range("D50").Value = "'" + '5.37'
$ws.Range("E50").Value = '  +13.45%  '
$ws.Range("D51").Value = "'" + '78.60'
$ws.Range("E51").Value = '  +3.34%  '
